# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sat Apr 20 03:09:45 UTC 2024 with GitHub Actions".
# Every touched cell is plain text in the source file (t="inlineStr", no
# explicit style), so each write goes through a text-number-format +
# Normal-style reset to stop Excel from auto-coercing numeric-looking
# strings (e.g. "556.52") into actual numbers / picking up a new style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '63.830.89'
Set-TextValue "E2" '  +4.23%  '
Set-TextValue "D3" '3.052.79'
Set-TextValue "E3" '  +3.66%  '
Set-TextValue "E4" '  +0.13%  '
Set-TextValue "D5" '556.52'
Set-TextValue "E5" '  +3.59%  '
Set-TextValue "D6" '142.34'
Set-TextValue "E6" '  +6.66%  '
Set-TextValue "E7" '  +0.26%  '
Set-TextValue "D8" '3.048.61'
Set-TextValue "E8" '  +3.41%  '
Set-TextValue "D9" '0.506'
Set-TextValue "E9" '  +5.55%  '
Set-TextValue "D10" '0.155'
Set-TextValue "E10" '  +7.57%  '
Set-TextValue "D11" '6.03'
Set-TextValue "E11" '  -7.78%  '
Set-TextValue "D12" '0.478'
Set-TextValue "E12" '  +9.12%  '
Set-TextValue "E13" '  +7.32%  '
Set-TextValue "D14" '34.91'
Set-TextValue "E14" '  +5.80%  '
Set-TextValue "D15" '3.553.23'
Set-TextValue "E15" '  +4.76%  '
Set-TextValue "D16" '63.877.41'
Set-TextValue "E16" '  +4.50%  '
Set-TextValue "E17" '  +3.63%  '
Set-TextValue "D18" '3.056.89'
Set-TextValue "E18" '  +4.29%  '
Set-TextValue "D19" '6.71'
Set-TextValue "E19" '  +3.61%  '
Set-TextValue "D20" '474.54'
Set-TextValue "E20" '  +3.77%  '
Set-TextValue "D21" '13.97'
Set-TextValue "E21" '  +6.37%  '
Set-TextValue "D22" '0.676'
Set-TextValue "E22" '  +5.70%  '
Set-TextValue "D23" '7.54'
Set-TextValue "E23" '  +8.28%  '
Set-TextValue "D24" '14.24'
Set-TextValue "E24" '  +16.56%  '
Set-TextValue "D25" '81.44'
Set-TextValue "E25" '  +4.38%  '
Set-TextValue "D26" '0.999'
Set-TextValue "E26" '  +0.00%  '
Set-TextValue "D27" '2.80'
Set-TextValue "E27" '  +4.91%  '
Set-TextValue "D28" '7.91'
Set-TextValue "E28" '  +7.05%  '
Set-TextValue "D29" '2.02'
Set-TextValue "E29" '  +4.32%  '
Set-TextValue "D30" '1.00'
Set-TextValue "E30" '  +0.26%  '
Set-TextValue "D31" '26.23'
Set-TextValue "E31" '  +5.74%  '
Set-TextValue "E32" '  +3.20%  '
Set-TextValue "D33" '2.43'
Set-TextValue "E33" '  +7.02%  '
Set-TextValue "D34" '5.58'
Set-TextValue "E34" '  +2.92%  '
Set-TextValue "D35" '6.19'
Set-TextValue "E35" '  +8.53%  '
Set-TextValue "D36" '54.82'
Set-TextValue "E36" '  +2.09%  '
Set-TextValue "D37" '0.0405'
Set-TextValue "E37" '  +6.90%  '
Set-TextValue "D38" '443.23'
Set-TextValue "E38" '  +0.87%  '
Set-TextValue "D39" '0.0805'
Set-TextValue "E39" '  +2.46%  '
Set-TextValue "D40" '2.81'
Set-TextValue "E40" '  +18.06%  '
Set-TextValue "D41" '2.964.70'
Set-TextValue "E41" '  +2.76%  '
Set-TextValue "D42" '8.19'
Set-TextValue "E42" '  +4.67%  '
Set-TextValue "D43" '0.113'
Set-TextValue "E43" '  +0.23%  '
Set-TextValue "D44" '27.50'
Set-TextValue "E44" '  +4.85%  '
Set-TextValue "D45" '0.259'
Set-TextValue "E45" '  +7.44%  '
Set-TextValue "B46" 'USDe'
Set-TextValue "C46" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D46" '1.00'
Set-TextValue "E46" '  -0.03%  '
Set-TextValue "B47" 'Fetch.AI'
Set-TextValue "C47" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D47" '2.14'
Set-TextValue "E47" '  +9.62%  '
Set-TextValue "E48" '  +5.39%  '
Set-TextValue "B49" 'Monero'
Set-TextValue "C49" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D49" '117.14'
Set-TextValue "E49" '  +3.75%  '
Set-TextValue "B50" 'PEPE'
Set-TextValue "C50" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D50" '0.0₃0512'
Set-TextValue "E50" '  +7.00%  '
Set-TextValue "D51" '2.06'
Set-TextValue "E51" '  +5.87%  '
